$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.34%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.69%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.700"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.32%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06224"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.51%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.729"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.12%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8526"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.42%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9102"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.51%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.79%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04879"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.58%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07085"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.97%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.11%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09052"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.73%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001533"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.21%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006160"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.18%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006001"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.12%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.452"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.11%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.172"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.14%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.51%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1309"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.87%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.103"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.21%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04239"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.47%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001220"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.20%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004080"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.01%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.35%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03920"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.94%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1111"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.30%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004111"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.57%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01387"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.23%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005114"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.38%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.00%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06956"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-47.34%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
